# Assignment 7 instructions: turn the plain "GitHub Project: " label into a
# bold label followed by the student's GitHub project URL, mirroring the
# existing "FlipGrid Link: (See Assignment #7 FlipGrid Wall)" paragraph
# right above it.

$d = $word.ActiveDocument

# Locate the "GitHub Project: " run and drop the trailing space (the bold
# label itself carries no trailing space once it becomes its own run).
$found = $d.Content.Find.Execute("GitHub Project: ", $false, $false, $false,
    $false, $false, $true, 1, $false, "GitHub Project:", 2)

$label = $d.Content
$null = $label.Find.Execute("GitHub Project:", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0)
$labelStart = $label.Start
$labelEnd = $label.End

# Append " " + the GitHub URL right after the label, still in the
# surrounding (non-bold) character formatting.
$tail = $d.Range($labelEnd, $labelEnd)
$tail.InsertAfter(" https://github.com/CSC4500/Jared-Heeringa")

# Now make just the label text bold, leaving the paragraph mark and the
# newly-appended space/URL untouched.
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Bold = 1
